$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the value that was in B4 (Transaction # amount column)
$ws.Range("B4").ClearContents()

# Move the active selection down to B5 (current working cell indicator)
$ws.Range("B5").Select()
